$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-7) contents but keep header row/style; extend through new row 13
$ws.Range("A2:T13").ClearContents()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2420556666666667
$ws.Range("H2").Value = 0.726167
$ws.Range("I2").Value = 0.5314769098578004
$ws.Range("J2").Value = 0.5314769098578004
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.118034
$ws.Range("N2").Value = 3.354102
$ws.Range("O2").Value = 0.01817106018145251
$ws.Range("P2").Value = 0.01817106018145251
$ws.Range("Q2").Value = 0.270626465226
$ws.Range("R2").Value = 2.435638187034
$ws.Range("S2").Value = 0.0096574989140785
$ws.Range("T2").Value = 0.009657498914078503

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2420556666666667
$ws.Range("H3").Value = 0.726167
$ws.Range("I3").Value = 0.5314769098578004
$ws.Range("J3").Value = 0.5314769098578004
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 59.24481466666666
$ws.Range("N3").Value = 177.734444
$ws.Range("O3").Value = 0.962887615892719
$ws.Range("P3").Value = 0.9628876158927191
$ws.Range("Q3").Value = 14.34054311068311
$ws.Range("R3").Value = 129.064887996148
$ws.Range("S3").Value = 0.5117525346350069
$ws.Range("T3").Value = 0.511752534635007

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2420556666666667
$ws.Range("H4").Value = 0.726167
$ws.Range("I4").Value = 0.5314769098578004
$ws.Range("J4").Value = 0.5314769098578004
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.121724666666667
$ws.Range("N4").Value = 3.365174
$ws.Range("O4").Value = 0.01823104344324033
$ws.Range("P4").Value = 0.01823104344324033
$ws.Range("Q4").Value = 0.2715198120064445
$ws.Range("R4").Value = 2.443678308058
$ws.Range("S4").Value = 0.009689378632696683
$ws.Range("T4").Value = 0.009689378632696683

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2420556666666667
$ws.Range("H5").Value = 0.726167
$ws.Range("I5").Value = 0.5314769098578004
$ws.Range("J5").Value = 0.5314769098578004
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04370233333333334
$ws.Range("N5").Value = 0.131107
$ws.Range("O5").Value = 0.0007102804825880949
$ws.Range("P5").Value = 0.0007102804825880949
$ws.Range("Q5").Value = 0.01057839742988889
$ws.Range("R5").Value = 0.095205576869
$ws.Range("S5").Value = 0.0003774976760182279
$ws.Range("T5").Value = 0.0003774976760182279

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna3"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2054156666666667
$ws.Range("H6").Value = 0.616247
$ws.Range("I6").Value = 0.4510271759376837
$ws.Range("J6").Value = 0.4510271759376837
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.118034
$ws.Range("N6").Value = 3.354102
$ws.Range("O6").Value = 0.01817106018145251
$ws.Range("P6").Value = 0.01817106018145251
$ws.Range("Q6").Value = 0.229661699466
$ws.Range("R6").Value = 2.066955295194
$ws.Range("S6").Value = 0.008195641957434219
$ws.Range("T6").Value = 0.00819564195743422

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna3"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2054156666666667
$ws.Range("H7").Value = 0.616247
$ws.Range("I7").Value = 0.4510271759376837
$ws.Range("J7").Value = 0.4510271759376837
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 59.24481466666666
$ws.Range("N7").Value = 177.734444
$ws.Range("O7").Value = 0.962887615892719
$ws.Range("P7").Value = 0.9628876158927191
$ws.Range("Q7").Value = 12.16981310129644
$ws.Range("R7").Value = 109.528317911668
$ws.Range("S7").Value = 0.4342884821414622
$ws.Range("T7").Value = 0.4342884821414622

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna3"
$ws.Range("C8").Value = "Epha3"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2054156666666667
$ws.Range("H8").Value = 0.616247
$ws.Range("I8").Value = 0.4510271759376837
$ws.Range("J8").Value = 0.4510271759376837
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.121724666666667
$ws.Range("N8").Value = 3.365174
$ws.Range("O8").Value = 0.01823104344324033
$ws.Range("P8").Value = 0.01823104344324033
$ws.Range("Q8").Value = 0.2304198202197778
$ws.Range("R8").Value = 2.073778381978
$ws.Range("S8").Value = 0.00822269603860191
$ws.Range("T8").Value = 0.00822269603860191

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna3"
$ws.Range("C9").Value = "Epha3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2054156666666667
$ws.Range("H9").Value = 0.616247
$ws.Range("I9").Value = 0.4510271759376837
$ws.Range("J9").Value = 0.4510271759376837
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.04370233333333334
$ws.Range("N9").Value = 0.131107
$ws.Range("O9").Value = 0.0007102804825880949
$ws.Range("P9").Value = 0.0007102804825880949
$ws.Range("Q9").Value = 0.008977143936555555
$ws.Range("R9").Value = 0.080794295429
$ws.Range("S9").Value = 0.0003203558001853635
$ws.Range("T9").Value = 0.0003203558001853635

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efna3"
$ws.Range("C10").Value = "Epha3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.007968333333333332
$ws.Range("H10").Value = 0.023905
$ws.Range("I10").Value = 0.01749591420451593
$ws.Range("J10").Value = 0.01749591420451593
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.118034
$ws.Range("N10").Value = 3.354102
$ws.Range("O10").Value = 0.01817106018145251
$ws.Range("P10").Value = 0.01817106018145251
$ws.Range("Q10").Value = 0.00890886759
$ws.Range("R10").Value = 0.08017980831
$ws.Range("S10").Value = 0.0003179193099397887
$ws.Range("T10").Value = 0.0003179193099397887

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Efna3"
$ws.Range("C11").Value = "Epha3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.007968333333333332
$ws.Range("H11").Value = 0.023905
$ws.Range("I11").Value = 0.01749591420451593
$ws.Range("J11").Value = 0.01749591420451593
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 59.24481466666666
$ws.Range("N11").Value = 177.734444
$ws.Range("O11").Value = 0.962887615892719
$ws.Range("P11").Value = 0.9628876158927191
$ws.Range("Q11").Value = 0.4720824315355555
$ws.Range("R11").Value = 4.24874188382
$ws.Range("S11").Value = 0.0168465991162499
$ws.Range("T11").Value = 0.0168465991162499

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Efna3"
$ws.Range("C12").Value = "Epha3"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.007968333333333332
$ws.Range("H12").Value = 0.023905
$ws.Range("I12").Value = 0.01749591420451593
$ws.Range("J12").Value = 0.01749591420451593
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.121724666666667
$ws.Range("N12").Value = 3.365174
$ws.Range("O12").Value = 0.01823104344324033
$ws.Range("P12").Value = 0.01823104344324033
$ws.Range("Q12").Value = 0.008938276052222222
$ws.Range("R12").Value = 0.08044448447
$ws.Range("S12").Value = 0.0003189687719417354
$ws.Range("T12").Value = 0.0003189687719417354

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Efna3"
$ws.Range("C13").Value = "Epha3"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.007968333333333332
$ws.Range("H13").Value = 0.023905
$ws.Range("I13").Value = 0.01749591420451593
$ws.Range("J13").Value = 0.01749591420451593
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.04370233333333334
$ws.Range("N13").Value = 0.131107
$ws.Range("O13").Value = 0.0007102804825880949
$ws.Range("P13").Value = 0.0007102804825880949
$ws.Range("Q13").Value = 0.0003482347594444445
$ws.Range("R13").Value = 0.003134112835
$ws.Range("S13").Value = 0.00001242700638450348
$ws.Range("T13").Value = 0.00001242700638450348

